# Auto-generated files on 2026-02-21
# Update the stock name cells in columns A-C for rows 2-21 on Sheet1
# to reflect the refreshed "Hot Stock Top 20" ranking lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @("天奇股份", "电广传媒", "风语筑")
    3  = @("电广传媒", "天奇股份", "博纳影业")
    4  = @("华胜天成", "华胜天成", "华胜天成")
    5  = @("网宿科技", "网宿科技", "天奇股份")
    6  = @("三花智控", "东方财富", "利欧股份")
    7  = @("风语筑", "三花智控", "汉缆股份")
    8  = @("利欧股份", "贵州茅台", "光线传媒")
    9  = @("汉缆股份", "利亚德", "三花智控")
    10 = @("光线传媒", "中国中铁", "掌阅科技")
    11 = @("大位科技", "利欧股份", "万向钱潮")
    12 = @("博纳影业", "风语筑", "嘉美包装")
    13 = @("首都在线", "百达精工", "深科技")
    14 = @("深科技", "汉缆股份", "协鑫集成")
    15 = @("万向钱潮", "光线传媒", "大位科技")
    16 = @("百达精工", "捷成股份", "二六三")
    17 = @("掌阅科技", "深科技", "紫金矿业")
    18 = @("二六三", "蓝色光标", "巨力索具")
    19 = @("东方财富", "万向钱潮", "百达精工")
    20 = @("贵州茅台", "掌阅科技", "卧龙电驱")
    21 = @("利亚德", "博纳影业", "国安股份")
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
}
